# Commit: "add stunting OR given diarrhoea"
#
# 1. Fix casing typo: "Incidence Diarrhoea" -> "Incidence diarrhoea"
# 2. Insert a new worksheet "OR stunting diarrhoea" right after the
#    "RR diarrhoea" sheet (and before "birth distribution"), holding a
#    single odds-ratio row (age-band headers + 1.04 across the board).

$wb = $excel.ActiveWorkbook

# --- 1) rename sheet -------------------------------------------------
$wb.Worksheets.Item("Incidence Diarrhoea").Name = "Incidence diarrhoea"

# --- 2) insert the new sheet in the right spot ------------------------
$afterSheet = $wb.Worksheets.Item("RR diarrhoea")
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$ws.Name = "OR stunting diarrhoea"

# --- 3) populate it ----------------------------------------------------
$ws.Range("A1").Value = "<1 month"
$ws.Range("B1").Value = "1-5 months"
$ws.Range("C1").Value = "6-11 months"
$ws.Range("D1").Value = "12-23 months"
$ws.Range("E1").Value = "24-59 months"

$ws.Range("A2").Value = 1.04
$ws.Range("B2").Value = 1.04
$ws.Range("C2").Value = 1.04
$ws.Range("D2").Value = 1.04
$ws.Range("E2").Value = 1.04
